$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.085237698923785
$ws.Cells.Item(2, 4).Value = 1.080387921002492
$ws.Cells.Item(2, 5).Value = 1.098412106269536
$ws.Cells.Item(2, 6).Value = 1.104248630094268
$ws.Cells.Item(2, 9).Value = 1.062205505190174
$ws.Cells.Item(2, 10).Value = 1.090093550749978
$ws.Cells.Item(2, 11).Value = 1.08306148668861
$ws.Cells.Item(2, 12).Value = 1.101039386913871
$ws.Cells.Item(2, 13).Value = 1.106861271754477
$ws.Cells.Item(2, 14).Value = 1.091641608627522
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.086762539613611
$ws.Cells.Item(3, 4).Value = 1.081605162027854
$ws.Cells.Item(3, 5).Value = 1.09992388394374
$ws.Cells.Item(3, 6).Value = 1.105816679204168
$ws.Cells.Item(3, 9).Value = 1.062755597980151
$ws.Cells.Item(3, 10).Value = 1.091278884427754
$ws.Cells.Item(3, 11).Value = 1.084096007792069
$ws.Cells.Item(3, 12).Value = 1.102371113842721
$ws.Cells.Item(3, 13).Value = 1.108250208109878
$ws.Cells.Item(3, 14).Value = 1.092828625615079
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.087747598767711
$ws.Cells.Item(4, 4).Value = 1.082391209728601
$ws.Cells.Item(4, 5).Value = 1.100900824813542
$ws.Cells.Item(4, 6).Value = 1.106830126128008
$ws.Cells.Item(4, 9).Value = 1.063109310588026
$ws.Cells.Item(4, 10).Value = 1.092043796377119
$ws.Cells.Item(4, 11).Value = 1.084763252286093
$ws.Cells.Item(4, 12).Value = 1.103231030583713
$ws.Cells.Item(4, 13).Value = 1.109147246451572
$ws.Cells.Item(4, 14).Value = 1.093594623827149
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.088161339044188
$ws.Cells.Item(5, 4).Value = 1.082721289545602
$ws.Cells.Item(5, 5).Value = 1.101311231554062
$ws.Cells.Item(5, 6).Value = 1.107255902955811
$ws.Cells.Item(5, 9).Value = 1.063257479683908
$ws.Cells.Item(5, 10).Value = 1.092364874032315
$ws.Cells.Item(5, 11).Value = 1.085043250178556
$ws.Cells.Item(5, 12).Value = 1.103592115503629
$ws.Cells.Item(5, 13).Value = 1.109523962257133
$ws.Cells.Item(5, 14).Value = 1.093916157449443
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.088230785887426
$ws.Cells.Item(6, 4).Value = 1.082776689595057
$ws.Cells.Item(6, 5).Value = 1.101380123318447
$ws.Cells.Item(6, 6).Value = 1.107327376778868
$ws.Cells.Item(6, 9).Value = 1.0632823268412
$ws.Cells.Item(6, 10).Value = 1.092418755751136
$ws.Cells.Item(6, 11).Value = 1.085090233206481
$ws.Cells.Item(6, 12).Value = 1.103652718578068
$ws.Cells.Item(6, 13).Value = 1.109587191296958
$ws.Cells.Item(6, 14).Value = 1.093970115686485
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.087753128663257
$ws.Cells.Item(7, 4).Value = 1.082395621736297
$ws.Cells.Item(7, 5).Value = 1.100906309855076
$ws.Cells.Item(7, 6).Value = 1.106835816454907
$ws.Cells.Item(7, 9).Value = 1.063111292515951
$ws.Cells.Item(7, 10).Value = 1.092048088555379
$ws.Cells.Item(7, 11).Value = 1.084766995635591
$ws.Cells.Item(7, 12).Value = 1.103235857076057
$ws.Cells.Item(7, 13).Value = 1.109152281706181
$ws.Cells.Item(7, 14).Value = 1.093598922100794
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.085753365130636
$ws.Cells.Item(8, 4).Value = 1.080799625676689
$ws.Cells.Item(8, 5).Value = 1.098923288142249
$ws.Cells.Item(8, 6).Value = 1.104778809726809
$ws.Cells.Item(8, 9).Value = 1.062391876157886
$ws.Cells.Item(8, 10).Value = 1.090494573861748
$ws.Cells.Item(8, 11).Value = 1.083411558392811
$ws.Cells.Item(8, 12).Value = 1.101489826968603
$ws.Cells.Item(8, 13).Value = 1.107331024763541
$ws.Cells.Item(8, 14).Value = 1.092043201238107
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.082216800851703
$ws.Cells.Item(9, 4).Value = 1.07797485382716
$ws.Cells.Item(9, 5).Value = 1.095418815241004
$ws.Cells.Item(9, 6).Value = 1.101144678430984
$ws.Cells.Item(9, 9).Value = 1.061106917701758
$ws.Cells.Item(9, 10).Value = 1.087740887320037
$ws.Cells.Item(9, 11).Value = 1.081006321212879
$ws.Cells.Item(9, 12).Value = 1.098399012042137
$ws.Cells.Item(9, 13).Value = 1.104108430750929
$ws.Cells.Item(9, 14).Value = 1.089285604145656
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.079850032304455
$ws.Cells.Item(10, 4).Value = 1.076082960195333
$ws.Cells.Item(10, 5).Value = 1.093075223245282
$ws.Cells.Item(10, 6).Value = 1.098715109313725
$ws.Cells.Item(10, 9).Value = 1.060238475440874
$ws.Cells.Item(10, 10).Value = 1.085893815990888
$ws.Cells.Item(10, 11).Value = 1.079391205330557
$ws.Cells.Item(10, 12).Value = 1.096328570800395
$ws.Cells.Item(10, 13).Value = 1.101950646300909
$ws.Cells.Item(10, 14).Value = 1.087435909763357
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.078822934189897
$ws.Cells.Item(11, 4).Value = 1.075261602126025
$ws.Cells.Item(11, 5).Value = 1.092058592589822
$ws.Cells.Item(11, 6).Value = 1.097661355301779
$ws.Cells.Item(11, 9).Value = 1.059859585383999
$ws.Cells.Item(11, 10).Value = 1.085091250442665
$ws.Cells.Item(11, 11).Value = 1.078689008885816
$ws.Cells.Item(11, 12).Value = 1.095429603639667
$ws.Cells.Item(11, 13).Value = 1.101013976068518
$ws.Cells.Item(11, 14).Value = 1.08663220448
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.078441073146297
$ws.Cells.Item(12, 4).Value = 1.074956182259982
$ws.Cells.Item(12, 5).Value = 1.091680685282598
$ws.Cells.Item(12, 6).Value = 1.097269674217945
$ws.Cells.Item(12, 9).Value = 1.059718416763145
$ws.Cells.Item(12, 10).Value = 1.084792717836671
$ws.Cells.Item(12, 11).Value = 1.078427748433628
$ws.Cells.Item(12, 12).Value = 1.09509531043858
$ws.Cells.Item(12, 13).Value = 1.100665695575298
$ws.Cells.Item(12, 14).Value = 1.086333247923463
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.078522999680866
$ws.Cells.Item(13, 4).Value = 1.075021710946068
$ws.Cells.Item(13, 5).Value = 1.091761760743747
$ws.Cells.Item(13, 6).Value = 1.097353703505256
$ws.Cells.Item(13, 9).Value = 1.05974871752392
$ws.Cells.Item(13, 10).Value = 1.084856773405542
$ws.Cells.Item(13, 11).Value = 1.078483809408446
$ws.Cells.Item(13, 12).Value = 1.095167034648946
$ws.Cells.Item(13, 13).Value = 1.10074041932875
$ws.Cells.Item(13, 14).Value = 1.086397394458589
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.078791376655471
$ws.Cells.Item(14, 4).Value = 1.075236362816013
$ws.Cells.Item(14, 5).Value = 1.092027360528621
$ws.Cells.Item(14, 6).Value = 1.097628984381813
$ws.Cells.Item(14, 9).Value = 1.059847925177226
$ws.Cells.Item(14, 10).Value = 1.08506658233878
$ws.Cells.Item(14, 11).Value = 1.078667421895821
$ws.Cells.Item(14, 12).Value = 1.095401978593749
$ws.Cells.Item(14, 13).Value = 1.100985194477274
$ws.Cells.Item(14, 14).Value = 1.086607501344578
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.078956685752986
$ws.Cells.Item(15, 4).Value = 1.075368572857165
$ws.Cells.Item(15, 5).Value = 1.092190967231175
$ws.Cells.Item(15, 6).Value = 1.09779855796726
$ws.Cells.Item(15, 9).Value = 1.059908992925094
$ws.Cells.Item(15, 10).Value = 1.085195796129709
$ws.Cells.Item(15, 11).Value = 1.07878049397588
$ws.Cells.Item(15, 12).Value = 1.095546685136471
$ws.Cells.Item(15, 13).Value = 1.101135960633842
$ws.Cells.Item(15, 14).Value = 1.08673689863391
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.079918148728742
$ws.Cells.Item(16, 4).Value = 1.076137425008185
$ws.Cells.Item(16, 5).Value = 1.093142654143528
$ws.Cells.Item(16, 6).Value = 1.09878500615741
$ws.Cells.Item(16, 9).Value = 1.060263560765613
$ws.Cells.Item(16, 10).Value = 1.085947020580215
$ws.Cells.Item(16, 11).Value = 1.079437747380499
$ws.Cells.Item(16, 12).Value = 1.096388179911411
$ws.Cells.Item(16, 13).Value = 1.102012760083566
$ws.Cells.Item(16, 14).Value = 1.087489189909303
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.08052063397252
$ws.Cells.Item(17, 4).Value = 1.076619123565292
$ws.Cells.Item(17, 5).Value = 1.093739123287053
$ws.Cells.Item(17, 6).Value = 1.099403308092846
$ws.Cells.Item(17, 9).Value = 1.060485206292425
$ws.Cells.Item(17, 10).Value = 1.086417496449979
$ws.Cells.Item(17, 11).Value = 1.079849259880762
$ws.Cells.Item(17, 12).Value = 1.096915364910024
$ws.Cells.Item(17, 13).Value = 1.102562121924176
$ws.Cells.Item(17, 14).Value = 1.087960333908766
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.080871834898749
$ws.Cells.Item(18, 4).Value = 1.076899882675289
$ws.Cells.Item(18, 5).Value = 1.094086856739539
$ws.Cells.Item(18, 6).Value = 1.099763786426212
$ws.Cells.Item(18, 9).Value = 1.0606142137785
$ws.Cells.Item(18, 10).Value = 1.086691650177847
$ws.Cells.Item(18, 11).Value = 1.080089014535199
$ws.Cells.Item(18, 12).Value = 1.09722262694131
$ws.Cells.Item(18, 13).Value = 1.10288233065152
$ws.Cells.Item(18, 14).Value = 1.088234876966374
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.080991548567729
$ws.Cells.Item(19, 4).Value = 1.076995579291343
$ws.Cells.Item(19, 5).Value = 1.094205395124954
$ws.Cells.Item(19, 6).Value = 1.09988667222377
$ws.Cells.Item(19, 9).Value = 1.060658155540291
$ws.Cells.Item(19, 10).Value = 1.086785084475738
$ws.Cells.Item(19, 11).Value = 1.080170718465449
$ws.Cells.Item(19, 12).Value = 1.097327355564372
$ws.Cells.Item(19, 13).Value = 1.102991475691988
$ws.Cells.Item(19, 14).Value = 1.088328443951686
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.080456015681416
$ws.Cells.Item(20, 4).Value = 1.076567463383548
$ws.Cells.Item(20, 5).Value = 1.093675146156468
$ws.Cells.Item(20, 6).Value = 1.099336987459453
$ws.Cells.Item(20, 9).Value = 1.060461454257026
$ws.Cells.Item(20, 10).Value = 1.086367046502483
$ws.Cells.Item(20, 11).Value = 1.079805136834546
$ws.Cells.Item(20, 12).Value = 1.0968588274259
$ws.Cells.Item(20, 13).Value = 1.102503203906059
$ws.Cells.Item(20, 14).Value = 1.087909812316558
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.078712356095737
$ws.Cells.Item(21, 4).Value = 1.075173162381936
$ws.Cells.Item(21, 5).Value = 1.091949155963379
$ws.Cells.Item(21, 6).Value = 1.0975479285514
$ws.Cells.Item(21, 9).Value = 1.05981872295299
$ws.Cells.Item(21, 10).Value = 1.085004810604206
$ws.Cells.Item(21, 11).Value = 1.078613364613156
$ws.Cells.Item(21, 12).Value = 1.095332803925691
$ws.Cells.Item(21, 13).Value = 1.100913124289306
$ws.Cells.Item(21, 14).Value = 1.086545641887055
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.077614012018516
$ws.Cells.Item(22, 4).Value = 1.074294591938371
$ws.Cells.Item(22, 5).Value = 1.09086230158654
$ws.Cells.Item(22, 6).Value = 1.096421509830987
$ws.Cells.Item(22, 9).Value = 1.0594121109546
$ws.Cells.Item(22, 10).Value = 1.084145861708309
$ws.Cells.Item(22, 11).Value = 1.0778615391038
$ws.Cells.Item(22, 12).Value = 1.094371148334863
$ws.Cells.Item(22, 13).Value = 1.099911294067123
$ws.Cells.Item(22, 14).Value = 1.085685473185204
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.07819646127411
$ws.Cells.Item(23, 4).Value = 1.074760522930815
$ws.Cells.Item(23, 5).Value = 1.091438623631512
$ws.Cells.Item(23, 6).Value = 1.097018797313884
$ws.Cells.Item(23, 9).Value = 1.05962790210053
$ws.Cells.Item(23, 10).Value = 1.084601442243338
$ws.Cells.Item(23, 11).Value = 1.078260336282076
$ws.Cells.Item(23, 12).Value = 1.09488114997478
$ws.Cells.Item(23, 13).Value = 1.100442583657715
$ws.Cells.Item(23, 14).Value = 1.086141700696848
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.080485214581129
$ws.Cells.Item(24, 4).Value = 1.076590807037412
$ws.Cells.Item(24, 5).Value = 1.093704055214788
$ws.Cells.Item(24, 6).Value = 1.099366955414356
$ws.Cells.Item(24, 9).Value = 1.060472187628035
$ws.Cells.Item(24, 10).Value = 1.086389843487018
$ws.Cells.Item(24, 11).Value = 1.079825074985541
$ws.Cells.Item(24, 12).Value = 1.096884375010966
$ws.Cells.Item(24, 13).Value = 1.102529827116998
$ws.Cells.Item(24, 14).Value = 1.087932641675426
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.083132647106508
$ws.Cells.Item(25, 4).Value = 1.078706634124457
$ws.Cells.Item(25, 5).Value = 1.096326053967808
$ws.Cells.Item(25, 6).Value = 1.102085354287413
$ws.Cells.Item(25, 9).Value = 1.061441175769652
$ws.Cells.Item(25, 10).Value = 1.088454740647201
$ws.Cells.Item(25, 11).Value = 1.081630156475764
$ws.Cells.Item(25, 12).Value = 1.098399012042137
$ws.Cells.Item(25, 13).Value = 1.104943170743216
$ws.Cells.Item(25, 14).Value = 1.090000471226424
